$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two header labels in row 1 (shared-string text changes;
# all other cells keep their existing text, they'll just be re-indexed
# automatically as the shared-string table is rebuilt).
$ws.Range("A1").Value = "labelForm"
$ws.Range("B1").Value = "labelType"

# Selection moves back to the sheet's default cell (A1) instead of the
# previously-saved A4.
$ws.Range("A1").Select()
